$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.738.38"

$ws.Range("D3").Value = "1.906.32"
$ws.Range("E3").Value = "  +0.68%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.49"
$ws.Range("E5").Value = "  -1.01%  "

$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4949"
$ws.Range("E7").Value = "  +0.75%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2956"
$ws.Range("E8").Value = "  +0.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06737"
$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("D10").Value = "1.907.65"
$ws.Range("E10").Value = "  +0.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.02"
$ws.Range("E11").Value = "  -0.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07356"
$ws.Range("E12").Value = "  +1.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.158"
$ws.Range("E13").Value = "  +2.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.31"
$ws.Range("E14").Value = "  -2.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6690"
$ws.Range("E15").Value = "  -0.90%  "

$ws.Range("D16").Value = "30.684.70"
$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007924"
$ws.Range("E17").Value = "  -0.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.50"
$ws.Range("E18").Value = "  +2.97%  "

$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").Value = "2.154.07"
$ws.Range("E20").Value = "  +0.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.306"
$ws.Range("E21").Value = "  +10.55%  "

$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "195.22"
$ws.Range("E23").Value = "  +2.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.259"
$ws.Range("E24").Value = "  +2.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.622"
$ws.Range("E25").Value = "  +2.82%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.99"
$ws.Range("E26").Value = "  +3.60%  "

$ws.Range("E27").Value = "  -1.31%  "

$ws.Range("E28").Value = "  +2.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.470"
$ws.Range("E29").Value = "  +4.82%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.431"
$ws.Range("E30").Value = "  +3.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09147"
$ws.Range("E31").Value = "  +0.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.042"
$ws.Range("E32").Value = "  +1.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05256"
$ws.Range("E33").Value = "  +0.42%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7409"
$ws.Range("E34").Value = "  +0.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.108"
$ws.Range("E35").Value = "  +0.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.722"
$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01820"
$ws.Range("E37").Value = "  -0.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.714"
$ws.Range("E38").Value = "  +1.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9216"
$ws.Range("E39").Value = "  -1.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.074"
$ws.Range("E40").Value = "  -2.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.18"
$ws.Range("E41").Value = "  +28.77%  "

$ws.Range("E42").Value = "  +0.96%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.925"
$ws.Range("E43").Value = "  +3.40%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.72"
$ws.Range("E44").Value = "  +1.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9992"
$ws.Range("E45").Value = "  -0.17%  "

$ws.Range("E46").Value = "  +2.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.549"
$ws.Range("E47").Value = "  +0.32%  "

$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.40"
$ws.Range("E48").Value = "  +4.80%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.068"
$ws.Range("E49").Value = "  +3.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05858"
$ws.Range("E50").Value = "  -0.18%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3998"
$ws.Range("E51").Value = "  +1.50%  "
